# Separate treatment support as a program
# Insert a new row for "program_prop_treatment_support_improvement" just
# above the existing "program_prop_lowquality" row (row 38) on the
# "constants" sheet, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert a new row at position 38; Excel copies the formatting of the row
# above (row 37), which is exactly the style the new row should have.
$ws.Rows(38).Insert()

$ws.Range("A38").Value = "program_prop_treatment_support_improvement"
$ws.Range("B38").Value = 0.401
$ws.Range("C38").Value = "Proportional reduction in adverse outcomes from the treatment support intervention"

# Match the saved view position/selection recorded after the edit.
$ws.Application.Goto($ws.Range("A22"), $true) | Out-Null
$ws.Range("C32").Select() | Out-Null
